# edit.ps1 - applies the "added startup script and prayer on the end of flyer" change
$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Apply the flyer-base paragraph style to the (only, empty) paragraph
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$p.Style = "flyer-base"

# -----------------------------------------------------------------
# 2. Section columns: <w:cols w:num="4" w:space="706"/>
# -----------------------------------------------------------------
$tc = $d.PageSetup.TextColumns
[void]$tc.SetCount(4)
$tc.Spacing = 35.3   # points -> 706 twips

# -----------------------------------------------------------------
# 3. New paragraph style "Czgwna" (Część główna)
# -----------------------------------------------------------------
$czesc = $d.Styles.Add("Czgwna", 1)
$czesc.NameLocal = "Część główna"
$czesc.ParagraphFormat.SpaceAfter = 0
$czesc.ParagraphFormat.LineSpacingRule = 0
$czescFont = $czesc.Font
$czescFont.Name = "Georgia"
$czescFont.NameFarEast = "Arial Unicode MS"
$czescFont.NameBi = "Arial Unicode MS"
$czescFont.Size = 9
$czescFont.SizeBi = 9
$czescFont.Color = 0
$czescFont.LanguageID = "pl-PL"
$czescFont.LanguageIDFarEast = "pl-PL"

# -----------------------------------------------------------------
# 4. New paragraph style "prayer"
# -----------------------------------------------------------------
$prayer = $d.Styles.Add("prayer", 1)
$prayer.NameLocal = "prayer"
$prayer.BaseStyle = "flyer-base"
$prayer.NextParagraphStyle = "verse"
$prayer.QuickStyle = $true
$prayer.ParagraphFormat.SpaceBefore = 3      # points -> 60 twips
$prayer.ParagraphFormat.LeftIndent = 14.4    # points -> 288 twips
$prayer.ParagraphFormat.RightIndent = 7.2    # points -> 144 twips
$prayer.ParagraphFormat.Alignment = 1        # wdAlignParagraphCenter
$prayer.Font.Size = 9

# -----------------------------------------------------------------
# 5. New character style "prayerZnak" (prayer Znak), linked to "prayer"
# -----------------------------------------------------------------
$prayerZnak = $d.Styles.Add("prayerZnak", 2)
$prayerZnak.NameLocal = "prayer Znak"
$prayerZnak.BaseStyle = "flyer-baseZnak"
$prayerZnak.Font.Name = "Source Sans Pro"
$prayerZnak.Font.Size = 9

# Link paragraph <-> character styles together (both directions)
$prayer.LinkStyle = "prayerZnak"
$prayerZnak.LinkStyle = "prayer"

Write-Output "edit applied"
